$d = $word.ActiveDocument

# 1. Remove the leading "This week " from the introduction paragraph.
$d.Content.Find.Execute("This week I was working on a package", $true, $false, $false, $false, $false, $true, 1, $false, "I was working on a package", 2)

# 2. Clean up the document: re-inserting the document's own "clean" OOXML
#    (as produced by WordOpenXML serialization) removes stray w:proofErr
#    spell/grammar-check markers and merges the runs that they had split,
#    exactly mirroring what Word does when it re-saves a document after the
#    reviewer/author accepted-and-cleared those proofing marks.
$xml = $d.Content.WordOpenXML
$d.Content.InsertXML($xml)
